$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "62.916.30"
$ws.Range("E2").Value = "  -1.30%  "

# Row 3
$ws.Range("D3").Value = "3.221.46"
$ws.Range("E3").Value = "  -1.37%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "527.85"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.45%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.97"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.88%  "

# Row 7
$ws.Range("E7").Value = "  +1.78%  "

# Row 8
$ws.Range("E8").Value = "  -0.11%  "

# Row 9
$ws.Range("D9").Value = "3.220.31"
$ws.Range("E9").Value = "  -1.07%  "

# Row 10
$ws.Range("E10").Value = "  -0.66%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.09"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -7.11%  "

# Row 12
$ws.Range("E12").Value = "  +3.23%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000252"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.75%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.10"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.88%  "

# Row 15
$ws.Range("D15").Value = "3.729.51"
$ws.Range("E15").Value = "  -1.03%  "

# Row 16
$ws.Range("E16").Value = "  -3.99%  "

# Row 17
$ws.Range("D17").Value = "3.217.87"
$ws.Range("E17").Value = "  -1.06%  "

# Row 18
$ws.Range("D18").Value = "62.735.00"
$ws.Range("E18").Value = "  -1.04%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.15"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.49%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.02"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.07%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.965"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.66%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "365.88"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.00%  "

# Row 23
$ws.Range("E23").Value = "  +4.40%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.31"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.92%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.02"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.28%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.95"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +5.74%  "

# Row 27
$ws.Range("E27").Value = "  +2.94%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.64"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.60%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.29"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.91%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.15"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.77%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "28.44"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.37%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "635.50"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.05%  "

# Row 33
$ws.Range("E33").Value = "  -2.79%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.25"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.96%  "

# Row 35
$ws.Range("E35").Value = "  +3.92%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "56.93"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.67%  "

# Row 37
$ws.Range("E37").Value = "  +0.05%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.67"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +4.10%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.375"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.82%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.24%  "

# Row 41
$ws.Range("D41").Value = "0.0₃0705"
$ws.Range("E41").Value = "  +12.77%  "

# Row 42
$ws.Range("E42").Value = "  +1.25%  "

# Row 43
$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.54"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +12.75%  "

# Row 44
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.867.16"
$ws.Range("E44").Value = "  +2.26%  "

# Row 45
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.93"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +11.58%  "

# Row 46
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.69"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +5.17%  "

# Row 47
$ws.Range("E47").Value = "  +3.69%  "

# Row 48
$ws.Range("E48").Value = "  -2.63%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.00"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +10.04%  "

# Row 50
$ws.Range("E50").Value = "  +2.12%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "134.52"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.21%  "
